$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 808, shifting rows 808:849 down to 809:850
$ws.Rows.Item(808).Insert()

# Populate the newly inserted row 808 with the new data point.
# Column A stores dates as plain text (matches the rest of the sheet), so
# force Text formatting before assigning, otherwise Excel auto-converts the
# "yyyy/mm/dd" string into a date serial number. ClearFormats afterwards so
# the cell keeps the plain/default style used by every other data row.
$ws.Cells.Item(808, 1).NumberFormat = "@"
$ws.Cells.Item(808, 1).Value = "2026/02/14"
$ws.Cells.Item(808, 1).ClearFormats()
$ws.Cells.Item(808, 2).Value = "土"
$ws.Cells.Item(808, 3).Value = 14
$ws.Cells.Item(808, 4).Value = 201
